$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.156.26', '  -1.82%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.831.28', '  -3.31%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.0000', '  +0.02%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '231.27', '  -3.16%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9997', '  -0.01%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4650', '  -3.90%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2694', '  -6.64%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06272', '  -4.38%  '),
    @(10, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.841.23', '  -2.59%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07379', '  -1.11%  '),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '16.05', '  -5.03%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.891', '  -4.42%  '),
    @(14, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '83.21', '  -5.65%  '),
    @(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6202', '  -7.57%  '),
    @(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.094.52', '  -1.94%  '),
    @(17, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.001', '  +0.09%  '),
    @(18, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '227.64', '  -2.85%  '),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007282', '  -4.05%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.38', '  -6.71%  '),
    @(21, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  +0.02%  '),
    @(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.068.17', '  -2.84%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.838', '  -8.52%  '),
    @(24, 'BitDAO', 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit', '0.3896', '  +7.59%  '),
    @(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.850', '  -5.67%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '164.69', '  -3.30%  '),
    @(27, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.093', '  -3.16%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.68', '  -6.18%  '),
    @(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.841', '  -6.32%  '),
    @(30, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1011', '  -2.03%  '),
    @(31, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.367', '  -2.43%  '),
    @(32, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.049', '  -7.00%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.760', '  -6.93%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04796', '  -5.62%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.126', '  -7.37%  '),
    @(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6994', '  -7.15%  '),
    @(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.687', '  -0.97%  '),
    @(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01814', '  -3.96%  '),
    @(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.610', '  -1.34%  '),
    @(40, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8936', '  -3.12%  '),
    @(41, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.925', '  -7.06%  '),
    @(42, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.0000', '  -0.33%  '),
    @(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '102.87', '  -4.03%  '),
    @(44, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.488', '  -2.85%  '),
    @(45, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3995', '  -7.23%  '),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.938', '  -6.72%  '),
    @(47, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1192', '  -6.98%  '),
    @(48, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '59.57', '  -7.53%  '),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '8.438', '  -6.56%  '),
    @(50, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05522', '  -2.72%  '),
    @(51, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '32.52', '  -4.61%  '),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
}
